$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.742.83'
$ws.Range("E2").Value = '  -2.42%  '

$ws.Range("D3").Value = '1.796.54'
$ws.Range("E3").Value = '  -1.83%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.001'
$ws.Range("E4").Value = '  +0.13%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '308.04'
$ws.Range("E5").Value = '  -2.03%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.000'
$ws.Range("E6").Value = '  +0.03%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4570'
$ws.Range("E7").Value = '  +1.71%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3709'
$ws.Range("E8").Value = '  -2.01%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07250'
$ws.Range("E9").Value = '  -3.36%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.8543'
$ws.Range("E10").Value = '  -4.84%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '20.35'
$ws.Range("E11").Value = '  -3.43%  '

$ws.Range("D12").Value = '1.795.80'
$ws.Range("E12").Value = '  -1.64%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.301'
$ws.Range("E13").Value = '  -2.37%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.489'
$ws.Range("E14").Value = '  -4.55%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.07035'
$ws.Range("E15").Value = '  -1.23%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '90.37'
$ws.Range("E16").Value = '  -4.46%  '

$ws.Range("E17").Value = '  +0.08%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008621'
$ws.Range("E18").Value = '  -2.46%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.9999'
$ws.Range("E19").Value = '  +0.04%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.61'
$ws.Range("E20").Value = '  -4.39%  '

$ws.Range("D21").Value = '26.763.78'
$ws.Range("E21").Value = '  -2.33%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.283'
$ws.Range("E22").Value = '  -0.52%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.60'
$ws.Range("E23").Value = '  -3.59%  '

$ws.Range("D24").Value = '2.018.84'
$ws.Range("E24").Value = '  -1.93%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.906'
$ws.Range("E25").Value = '  -4.94%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '149.33'
$ws.Range("E26").Value = '  -1.74%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.151'
$ws.Range("E27").Value = '  -14.58%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.13'
$ws.Range("E28").Value = '  -2.99%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.199'
$ws.Range("E29").Value = '  -3.71%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '114.03'

$ws.Range("E31").Value = '  -0.23%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.7544'
$ws.Range("E32").Value = '  -3.84%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.157'
$ws.Range("E33").Value = '  -3.93%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.433'
$ws.Range("E34").Value = '  -3.51%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.878'
$ws.Range("E35").Value = '  -0.34%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.9995'
$ws.Range("E36").Value = '  -0.03%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.113'
$ws.Range("E37").Value = '  -0.10%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01938'
$ws.Range("E38").Value = '  -2.74%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.05204'
$ws.Range("E39").Value = '  -2.78%  '

$ws.Range("B40").Value = 'MXToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.899'
$ws.Range("E40").Value = '  +1.44%  '

$ws.Range("B41").Value = 'RenderToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.365'
$ws.Range("E41").Value = '  +2.53%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '7.132'
$ws.Range("E42").Value = '  -4.10%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.5225'
$ws.Range("E43").Value = '  -2.55%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.1644'
$ws.Range("E44").Value = '  -5.24%  '

$ws.Range("B45").Value = 'Aptos'
$ws.Range("C45").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '8.482'
$ws.Range("E45").Value = '  -4.10%  '

$ws.Range("B46").Value = 'Decentraland'
$ws.Range("C46").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.5001'
$ws.Range("E46").Value = '  -3.75%  '

$ws.Range("E47").Value = '  -5.44%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '104.05'
$ws.Range("E48").Value = '  -2.60%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.9994'

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.645'
$ws.Range("E50").Value = '  -3.97%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.06292'
$ws.Range("E51").Value = '  -1.66%  '
